# The "Förändrad" column (C) stores a date serial number that is bumped by
# one day (46060 -> 46061) for every data row (rows 2 through 161).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C161").Value = 46061
